$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = -8
$ws.Range("F3").Value = -6
$ws.Range("F4").Value = -1
$ws.Range("F8").Value = -13
$ws.Range("F9").Value = 0
$ws.Range("F10").Value = -1
